$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q factor data run for sg_rr_52_025 2023-12-11 18-30-06.csv data
# Append a new data row (row 82) to the table below the last existing
# entry (row 81), repeating the filename/scan-config values that are
# shared by the other rows in this block and filling in the Q-factor
# results computed for this run.
$ws.Range("A82").Value = "sg_rr_52_025 2023-12-11 18-30-06.csv"
$ws.Range("B82").Value = 0.01
$ws.Range("C82").Value = 1000
$ws.Range("D82").Value = 5001
$ws.Range("E82").Value = 1530
$ws.Range("F82").Value = 1570
$ws.Range("G82").Value = 0.01
$ws.Range("H82").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = 1.8875
$ws.Range("K82").Value = 0.0078430324425366096
$ws.Range("L82").Value = "yes"
$ws.Range("M82").Value = 0.151131672031217
$ws.Range("N82").Value = 0.0060602112450080696
$ws.Range("O82").Value = 10509.3396239214
$ws.Range("P82").Value = 323.17436740347603
$ws.Range("Q82").Value = 895614837.28920305
$ws.Range("R82").Value = 82773263.426702306
$ws.Range("S82").Value = 52
$ws.Range("T82").Value = 0.1

# Scroll the view down to reveal the newly added row, and leave the
# selection on A57 (matches the saved view state in the workbook).
$win = $excel.ActiveWindow
$win.ScrollRow = 71
$win.ScrollColumn = 1
$ws.Range("A57").Select()
